# Generate Report for Handback
# - Flip status text from "Ready for handoff" to "Handed back: in sync with en-US"
# - Widen a few columns that now need to show full file names
# - Fill in "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#   for each localized-language report sheet, linking the target file the same
#   way the source file is already linked.

$wb = $excel.ActiveWorkbook

$urlFor82304 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/489a8dd4042eb2fdff743600fb288cee8ee7fe4c/e2e/82304c21-bfe5-443d-8050-97c14a7bf7ad.md"
$urlForD8a125 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/489a8dd4042eb2fdff743600fb288cee8ee7fe4c/e2e/d8a125cf-0be1-4789-a552-c4854026fd03.md"

# 1) Status text: handback is now complete and in sync with en-US.
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US")
}

# 2) Widen columns that will hold the (long) target/handback file names.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

# 3) zh-cn report rows: record the handed-back target file + datetime.
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $urlFor82304, "", "", "82304c21-bfe5-443d-8050-97c14a7bf7ad.md")
$zhcn.Range("J2").Value = $zhcn.Range("G2").Value()
$zhcn.Range("K2").Value = "2016-09-07 03:18:19"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $urlForD8a125, "", "", "d8a125cf-0be1-4789-a552-c4854026fd03.md")
$zhcn.Range("J3").Value = $zhcn.Range("G3").Value()
$zhcn.Range("K3").Value = "2016-09-07 03:18:19"

# 4) de-de report rows: record the handed-back target file + datetime.
$dede.Hyperlinks.Add($dede.Range("I2"), $urlFor82304, "", "", "82304c21-bfe5-443d-8050-97c14a7bf7ad.md")
$dede.Range("J2").Value = $dede.Range("G2").Value()
$dede.Range("K2").Value = "2016-09-07 03:18:27"

$dede.Hyperlinks.Add($dede.Range("I3"), $urlForD8a125, "", "", "d8a125cf-0be1-4789-a552-c4854026fd03.md")
$dede.Range("J3").Value = $dede.Range("G3").Value()
$dede.Range("K3").Value = "2016-09-07 03:18:27"
